$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 to shift existing data down by one row
$ws.Rows.Item(2).Insert()
$ws.Range("B2:E2").ClearFormats()

# Column A holds forecast dates and carries the centered date-format style
# used by the rest of the date column (A3:A19)
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").Borders.LineStyle = 1

# Full target dataset for rows 2..19 (columns A-E)
# A = date_of_forecast, B = y_0, C = y_0_forecast, D = y_1, E = y_1_forecast
$data = @(
    @(39400, 2007, 3.14593994906931,    2008, $null),
    @(39765, 2008, 1.769627576887389,   2009, $null),
    @(40130, 2009, -4.774178217057756,  2010, $null),
    @(40494, 2010, 1.97975191822708,    2011, 2.743551941645217),
    @(40862, 2011, 3.452886745653183,   2012, 1.794132456841213),
    @(41228, 2012, 1.239479831392853,   2013, 1.158731032337301),
    @(41592, 2013, 0.2379616621361214,  2014, 1.348985046565354),
    @(41957, 2014, 1.51977456621637,    2015, 0.7739869831243862),
    @(42321, 2015, 1.470039379455756,   2016, 1.319057785023592),
    @(42689, 2016, 1.638797242243251,   2017, 1.236938064849924),
    @(43053, 2017, 2.161565493242668,   2018, 2.320541194291881),
    @(43418, 2018, 2.214251681313772,   2019, 0.6374750548026054),
    @(43783, 2019, 0.6066442151010376,  2020, 0.57214245765278),
    @(44159, 2020, -4.207901339433196,  2021, -0.2586890779524231),
    @(44525, 2021, 1.099928004397532,   2022, 1.794400784768979),
    @(44890, 2022, 2.310042359896225,   2023, 0.8232644777432796),
    @(45254, 2023, 0.0464415346324687,  2024, 0.4457784880425031),
    @(45618, 2024, -0.3101476031197148, 2025, 0.5215192790195111)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    if ($row[4] -ne $null) {
        $ws.Cells.Item($r, 5).Value = $row[4]
    } else {
        $ws.Cells.Item($r, 5).ClearContents()
    }
    $r = $r + 1
}
